# August 2020 Measles Caseload Update
# HR measles data ingest - fill in the new "August" reporting row (row 40)
# and append the corresponding monthly-summary row (row 66).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New reporting period row (row 40): reporting start/end dates, totals,
#     confirmed cases, deaths, month label, and the two delta formulas that
#     mirror the pattern used by every other row in the table.
$ws.Range("A40").Value = 43831
$ws.Range("B40").Value = 44052
$ws.Range("C40").Value = 69249
$ws.Range("D40").Value = 1317
$ws.Range("E40").Value = 987
$ws.Range("F40").Value = "August"
$ws.Range("G40").Formula = "=SUM(C40-C39)"
$ws.Range("H40").Formula = "=SUM(E40-E39)"

# --- New monthly-summary row appended at the bottom of the sheet (row 66),
#     matching the August figures just computed above.
$ws.Range("A66").Value = "August"
$ws.Range("B66").Value = 623
$ws.Range("C66").Value = 19

# --- View-state bookkeeping: scroll/select so the new rows are in view,
#     matching the author's on-screen position when the file was saved.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1
$ws.Range("D63").Select()
